# Optimized the code in finding violations and plate numbers
#
# - SWERVING sheet gains the row that used to live on DRUNK DRIVING
#   (plate "xyz" [corrected from the typo "zxy"], class truck, color black,
#   date December 12, 2013, time 12:00nn)
# - DRUNK DRIVING sheet gets a brand new row 2 for a newly found violation
#   (plate def456, class sedan, color red, date January 10, 2016, time 5:26pm)
# - Cursor/selection positions shift around on several sheets as a result
#   of the edits.

$wb = $excel.ActiveWorkbook

$speeding       = $wb.Worksheets.Item("SPEEDING")
$swerving       = $wb.Worksheets.Item("SWERVING")
$drunkDriving   = $wb.Worksheets.Item("DRUNK DRIVING")
$counterflowing = $wb.Worksheets.Item("COUNTERFLOWING")

# --- SWERVING: add new row 2 (data moved over from DRUNK DRIVING, with the
#     plate number typo "zxy" corrected to "xyz") ---
# Columns are: A=PLATE NUMBER, B=CLASS, C=COLOR, D=DATE, E=TIME
$swerving.Range("A2").Value = "xyz"
$swerving.Range("B2").Value = "truck"
$swerving.Range("C2").Value = "black"
$swerving.Range("D2").Value = "December 12, 2013"
$swerving.Range("E2").Value = "12:00nn"

# --- DRUNK DRIVING: replace row 2 with a newly found violation ---
$drunkDriving.Range("A2").Value = "def456"
$drunkDriving.Range("B2").Value = "sedan"
$drunkDriving.Range("C2").Value = "red"
$drunkDriving.Range("D2").Value = "January 10, 2016"
$drunkDriving.Range("E2").Value = "5:26pm"

# --- Selection / active-cell changes left over from editing ---
$speeding.Activate()
$speeding.Range("A2").Select()

$swerving.Activate()
$swerving.Range("A2:E2").Select()

$counterflowing.Activate()
$counterflowing.Range("A2").Select()

# DRUNK DRIVING stays the active/tab-selected sheet, as in the original file.
$drunkDriving.Activate()
$drunkDriving.Range("D2").Select()
